$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$staffText = @'
อาจารย์และบุคลากร คณะศึกษาศาสตร์และนวัตกรรมการศึกษา
ผู้บริหารคณะศึกษาศาสตร์และนวัตกรรมการศึกษา
- ผศ.ดร.ลาวัณย์ ดุลยชาติ (คณบดี)
- ผศ.ดร.วทัญญู แก้วสุพรรณ (รองคณบดีงานบริหารและวางแผน)
- ผศ.นคินทร พัฒนชัย (รองคณบดีงานวิชาการและวิจัย)
- ดร.ปัญญา เถาว์ชาลี (รองคณบดีงานกิจการนักศึกษา)
- ผศ.ดร.ทิพย์อุบล ทิพเลิศ (ผู้ช่วยคณบดีงานโครงการจัดตั้งโรงเรียนสาธิต มหาวิทยาลัยกาฬสินธุ์)
- อาจารย์รศรงค์ พัฒนาอนุสรณ์ (ผู้ช่วยคณบดีงานศูนย์ความเป็นเลิศด้านนวัตกรรมการจัดการเรียนรู้)
- อาจารย์ธีรภัทร สินธุเดช (ผู้ช่วยคณบดีงานประกันคุณภาพการศึกษา)
- ผศ.วรนุช นิลเขต (ผู้ช่วยคณบดีงานบริหารโครงการครูรัก(ษ์)ถิ่น)
- อาจารย์กมลพัฒน์ ไชยสงคราม (ผู้ช่วยคณบดีงานวิเทศสัมพันธ์)
- อาจารย์ปรีชา ทับสมบัติ (ผู้ช่วยคณบดีงานกิจการพิเศษ)
- อาจารย์สมใจ ภูครองทุ่ง (ผู้ช่วยคณบดีงานฝึกประสบการณ์วิชาชีพครู)
- ว่าที่ร้อยตรีวรพัทธนันท์ ศรีสูงเนิน (ปฏิบัติหน้าที่หัวหน้าสำนักงานคณบดี)
อาจารย์ประจำหลักสูตรครุศาสตรบัณฑิต สาขาวิชานวัตกรรมการจัดการเรียนรู้ ระดับปริญญาตรี 12 วิชาเอก
-ผศ.ดร.สายหยุด ภูปุย (หัวหน้าสาขาวิชานวัตกรรมการจัดการเรียนรู้)
1.วิชาเอกการศึกษาปฐมวัย
- จำลองลักษณ์ เสียงสนั่น (หัวหน้าวิชาเอกการศึกษาปฐมวัย)
- ผศ.ดร.สายหยุด ภูปุย
- ผศ.วรนุช นิลเขต
- อาจารย์สุกานดา ภูจีระ
- อาจารย์ณัฐฐิยา ภูมิโยชน์
2.วิชาเอกการประถมศึกษา
- ผศ.ดร.ณิชาภาท์ กันขุนทศ (หัวหน้าวิชาเอกการประถมศึกษา)
- ดร.สุพจน์ ดวงเนตร
- อาจารย์เยาวเรศ รัตนธารทอง
- ผศ.ดร.นาตยา หกพันนา
- อาจารย์ศศิธร แสนพันดร
3.วิชาเอกภาษาไทย
- วัชรวร วงศ์กัณหา (หัวหน้าวิชาเอกภาษาไทย)
- อาจารย์ธีรภัทร สินธุเดช
- ผศ.ดร.อนุชา พิมศักดิ์
- อาจารย์นิวัฒน์ ชินเสริม
- อาจารย์จารญา อนันตะวัน
4.วิชาเอกภาษาอังกฤษ
- อาจารย์นฤตา หงษ์ษา (หัวหน้าวิชาเอกภาษาอังกฤษ)
- อาจารย์นันทณัฏฐ์ เวียงอินทร์
- ผศ.ศศิกร สุรมณี
- อาจารย์กมลพัฒน์ ไชยสงคราม
- อาจารย์ภัทราภรณ์ วาทะวัฒนะ
5.วิชาเอกคอมพิวเตอร์
- ว่าที่ ร.ต.สุรจักษ์ พิริยะเชิดชูชัย (หัวหน้าวิชาเอกคอมพิวเตอร์)
- ผศ.นคินทร พัฒนชัย
- ผศ.ดร.สวียา สุรมณี
- ผศ.ดร.ลาวัณย์ ดุลยชาติ
- ผศ.ดร.อัญญปารย์ ศิลปนิลมาลย์
- อาจารย์ธีรชาติ น้อยสมบัติ
- อาจารย์รศรงค์ พัฒนาอนุสรณ์
6.วิชาเอกวิทยาศาสตร์ทั่วไป
- ผศ.ตะวัน ทองสุข (หัวหน้าวิชาเอกวิทยาศาสตร์ทั่วไป)
- ผศ.ดร.ชุลิดา เหมตะศิลป์
- ผศ.ดร.วิศรุต พยุงเกียรติคุณ
- ผศ.อังคาร อินทนิล
7.วิชาเอกคณิตศาสตร์
- ผศ.ประภาพร หนองหารพิทักษ์ (หัวหน้าวิชาเอกคณิตศาสตร์)
- ผศ.ดร.ปนัดดา สังข์ศรีแก้ว
- ผศ.ดร.ปวีณา ขันธ์ศิลา
- ผศ.ดร.วรรณธิดา ยลวิลาศ
- อาจารย์สมใจ ภูครองทุ่ง
- ผศ.สุวรรณวัฒน์ เทียนยุทธกุล
8.วิชาเอกสังคมศึกษา
- ดร.นิติยา ค้อไผ่ (หัวหน้าวิชาเอกสังคมศึกษา)
- ดร.ปัญญา เถาว์ชาลี
- ผศ.ดร.แจ่มจันทร์ ณ กาฬสินธุ์
- อาจารย์อภิรดี ดอนอ่อนเบ้า
- อาจารย์ธนาคาร ผินสู่
9.วิชาเอกเกษตรศาสตร์
- ดร.สุภาพร พุ่มริ้ว (หัวหน้าวิชาเอกเกษตรศาสตร์)
- ผศ.ดร.กีรวิชญ์ เพชรจุล อาจารย์(สังกัดคณะเทคโนโลยีการเกษตร)
- ผศ.ดร.ชโลธร อัมพร อาจารย์(สังกัดคณะเทคโนโลยีการเกษตร)
- ดร.วรมัน ไม้เจริญ อาจารย์(สังกัดคณะเทคโนโลยีการเกษตร)
- ว่าที่ร้อยตรี ดร.ธนภูมิ บุญมี อาจารย์(สังกัดคณะเทคโนโลยีการเกษตร)
- ดร.ประสิทธิ์ ขุนสนิท อาจารย์(สังกัดคณะเทคโนโลยีการเกษตร)
10.วิชาเอกชีววิทยา
- ผศ.ดร.กชพรรณ วงค์เจริญ (หัวหน้าวิชาเอกชีววิทยา) อาจารย์(คณะวิทยาศาสตร์และเทคโนโลยีสุขภาพ)
- ผศ.ดร.ทิพย์อุบล ทิพเลิศ
- ผศ.ดร.ระพีพรรณ ประจันตะเสน อาจารย์(คณะวิทยาศาสตร์และเทคโนโลยีสุขภาพ)
11.วิชาเอกพลศึกษา
- อาจารย์ปรีชา ทับสมบัติ (หัวหน้าวิชาเอกพลศึกษา)
- อาจารย์ธนัญชัย เฉลิมสุข
- ผศ.ดร.วทัญญู แก้วสุพรรณ
- อาจารย์ภัทรนิณทร์ เหล่าแสงสา
12.วิชาเอกฟิสิกส์
- ผศ.ศิริพร จรรยา (หัวหน้าวิชาเอกฟิสิกส์)
- อาจารย์จตุพร คำสงค์ อาจารย์(คณะวิทยาศาสตร์และเทคโนโลยีสุขภาพ)
- อาจารย์เจษฎา ขจรฤทธิ์ อาจารย์(คณะวิทยาศาสตร์และเทคโนโลยีสุขภาพ)
อาจารย์ประจำหลักสูตรประกาศนียบัตรบัณฑิต สาขาวิชานวัตกรรมการจัดการเรียนรู้
- ดร.สุพจน์ ดวงเนตร (ประธานหลักสูตร)
- ผศ.ดร.แจ่มจันทร์ ณ กาฬสินธุ์
- ผศ.ดร.ณิชาภาท์ กันขุนทศ
- ผศ.ดร.คมสันทิ์ ขจรปัญญาไพศาล
- ผศ.ดร.ลาวัณย์ ดุลยชาติ
- ผศ.ดร.ศักดิ์สิทธิ์ ฤทธิลัน
- ผศ.ดร.สายหยุด ภูปุย
- ผศ.ดร.วรรณธิดา ยลวิลาศ
- ผศ.ดร.ปนัดดา สังข์ศรีแก้ว
- ดร.ปัญญา เถาว์ชาลี
- ผศ.ดร.ทิพย์อุบล ทิพเลิศ
- อาจารย์นันทนัฏฐ์ เวียงอินทร์
- ผศ.ศศิกร สุรมณี
- อาจารย์กมลพัฒน์ ไชยสงคราม
- ผศ.ดร.สวียา สุรมณี
- ผศ.ดร.อัญญปารย์ ศิลปนิลมาลย์
- ผศ.ดร.นาตยา หกพันนา
- ผศ.ดร.อมร มะลาศรี
- อาจารย์เยาวเรศ รัตนธารทอง
อาจารย์ประจำหลักสูตรประกาศนียบัตรบัณฑิตชั้นสูง สาขาวิชานวัตกรรมการจัดการเรียนรู้
- ผศ.ดร.ปนัดดา สังข์ศรีแก้ว (ประธานหลักสูตร)
- ผศ.ดร.ลาวัณย์ ดุลยชาติ
- ผศ.ดร.ปวีณา ขันธ์ศิลา
- ผศ.ดร.สายหยุด ภูปุย
- ผศ.ดร.แจ่มจันทร์ ณ กาฬสินธุ์
- ผศ.ดร.อนุชา พิมศักดิ์
- ผศ.ดร.คมสันทิ์ ขจรปัญญาไพศาล
- ผศ.ดร.ศักดิ์สิทธิ์ ฤทธิลัน
- ผศ.ดร.วรรณธิดา ยลวิลาศ
- ดร.สุพจน์ ดวงเนตร
อาจารย์ประจำหลักสูตรครุศาสตรมหาบัณฑิต สาขาวิชาการบริหารการศึกษา ระดับปริญญาโท
- ผศ.ดร.คมสันทิ์ ขจรปัญญาไพศาล (หัวหน้าสาขาวิชาการบริหารการศึกษา/ประธานหลักสูตร)
- ผศ.ดร.อมร มะลาศรี
- ผศ.ดร.ศักดิ์สิทธิ์ ฤทธิลัน
อาจารย์ประจำหลักสูตรครุศาสตรมหาบัณฑิต สาขาวิชานวัตกรรมการจัดการเรียนรู้ ระดับปริญญาโท
- ผศ.ดร.สายหยุด ภูปุย (ประธานหลักสูตร)
- ผศ.ดร.วรรณธิดา ยลวิลาศ
- ผศ.ดร.อนุชา พิมศักดิ์
- ผศ.ดร.อัญญปารย์ ศิลปนิลมาลย์
- ผศ.ดร.ปวีณา ขันธ์ศิลา
- ผศ.ดร.ชุลิดา เหมตะศิลป์
- ผศ.ดร.วิศรุต พยุงเกียรติคุณ
อาจารย์ประจำหลักสูตรครุศาสตรดุษฎีบัณฑิต สาขาวิชานวัตกรรมการจัดการเรียนรู้ ระดับปริญญาเอก
- ผศ.ดร.วิศรุต พยุงเกียรติคุณ (ประธานหลักสูตร)
- ผศ.ดร.สายหยุด ภูปุย
- ผศ.ดร.อนุชา พิมศักดิ์
- ผศ.ดร.อัญญปารย์ ศิลปนิลมาลย์
- ผศ.ดร.ปวีณา ขันธ์ศิลา
- ผศ.ดร.ชุลิดา เหมตะศิลป์
- ผศ.ดร.วรรณธิดา ยลวิลาศ
บุคลากรสายสนับสนุน สำนักงานคณบดี คณะศึกษาศาสตร์และนวัตกรรมการศึกษา
1.ว่าที่ร้อยตรีหญิงวรพัทธนันท์ ศรีสูงเนิน (ตำแหน่งนักวิชาการเงินและบัญชี ชำนาญการ ปฏิบัติหน้าที่หัวหน้าสำนักงาน และหัวหน้างานบริหารและวางแผน)
2.นางสาวภารดี จันทร์ลอย (ตำแหน่งเจ้าหน้าที่บริหารงานทั่วไป หัวหน้างานวิชาการและวิจัย รับผิดชอบงานประกันคุณภาพการศึกษา)
3.นายสัณห์ แทบพล (ตำแหน่งนักวิชาการศึกษา หัวหน้างานกิจการนักศึกษา รับผิดชอบงานกิจการนักศึกษา)
4.นางสาวปนัดดา ดวงเพชรแสง (ตำแหน่งเจ้าหน้าที่บริหารงานทั่วไป ชำนาญการ รับผิดชอบงานพัสดุและสินทรัพย์)
5.นางสาวสุนันทา จันมีวงษ์ (ตำแหน่งนักวิชาการศึกษา ชำนาญการ รับผิดชอบงานบัณฑิตศึกษาระดับปริญญาโท-เอก สาขาวิชานวัตกรรมการจัดการเรียนรู้)
6.นางสาวปิยนุช เกี้ยนมา (ตำแหน่งเจ้าหน้าที่บริหารงานทั่วไป รับผิดชอบงานวิชาการระดับปริญญาตรี)
7.นางธนิกานต์ วินิจ (ตำแหน่งเจ้าหน้าที่บริหารงานทั่วไป รับผิดชอบงานธุรการและสารบรรณ)
8.นางสาวนิตยา บุญนามน (ตำแหน่งเจ้าหน้าที่บริหารงานทั่วไป รับผิดชอบงานบัณฑิตศึกษาระดับปริญญาโท สาขาวิชาการบริหารการศึกษา)
9.นายอนุชิต คำหินกอง (ตำแหน่งนักวิชาการศึกษา รับผิดชอบงานฝึกประสบการณ์วิชาชีพครู)
10.นางสาววิรัชดา พรมคำบุตร (ตำแหน่งนักวิชาการศึกษา รับผิดชอบงานธุรการงานกิจการนักศึกษา)
11.นางสาวสุกัญญา ขะกิจ (ตำแหน่งนักวิชาการศึกษา รับผิดชอบงานประชาสัมพันธ์และสื่อสารองค์กร และเป็นผู้สร้าง AI Chatbot ที่คุณกำลังใช้อยู่นี้)
12.นางสาวอณัฐนงค์ พูลจรัส (ตำแหน่งเจ้าหน้าที่บริหารงานทั่วไป รับผิดชอบงานอาคารสถานที่และงานพัสดุโครงการครูรัก(ษ์)ถิ่น)
13.นายประเสริฐศักดิ์ วันนุกูล (ตำแหน่งนักวิชาการโสตทัศนศึกษา รับผิดชอบงานโสตทัศนูปกรณ์)
14.นางสาวศุภลักษณา ภูกาสอน (ตำแหน่งเจ้าหน้าที่บริหารงานทั่วไป รับผิดชอบงานบริหารและวางแผน)
15.นางสาวนาฎนภา นามเขต (ตำแหน่งนักวิชาการศึกษา รับผิดชอบงานโครงการครูรัก(ษ์)ถิ่น)
'@

$ws.Range("A9").Value2 = $staffText
$ws.Rows.Item(8).RowHeight = 15.4
$ws.Rows.Item(9).RowHeight = 16.5
$ws.Range("A9").Font.Bold = $false
$ws.Range("A9").Select() | Out-Null
